$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.844.74"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.817.36"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "663.32"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.18"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.816.69"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.21%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.77"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.461.87"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.799.33"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.830.18"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.81"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.18"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.21"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.32"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.75%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.74%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.29"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.86%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.968.70"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.84"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.32%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.39%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.62"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.66%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +15.03%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.773.63"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.48%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.92"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("B44").Value = "USDe"

$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Stacks"

$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.09"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +9.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.39"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.80"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.26%  "

$ws.Range("B48").Value = "ONDO"

$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.44"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.11%  "

$ws.Range("B49").Value = "OKB"

$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.91"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("B50").Value = "TheGraph"

$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.301"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000291"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.69%  "
